$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns A and B (header + data rows 1-19)
for ($r = 1; $r -le 19; $r++) {
    $valA = $ws.Cells.Item($r, 1).Value2
    $valB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $valB
    $ws.Cells.Item($r, 2).Value = $valA
}
